# Fix list level numbering: decrement the indent/outline level of every
# paragraph in the body placeholders on slides 1 and 2 by one level, so
# that top-level list items sit at the same level as top-level paragraphs
# (PowerPoint's TextRange.IndentLevel is 1-based, so "decrement the pPr
# lvl attribute by 1" means "decrement IndentLevel by 1" too).

$p = $ppt.ActivePresentation

for ($slideIdx = 1; $slideIdx -le $p.Slides.Count; $slideIdx++) {
    $s = $p.Slides.Item($slideIdx)
    for ($shpIdx = 1; $shpIdx -le $s.Shapes.Count; $shpIdx++) {
        $shp = $s.Shapes.Item($shpIdx)
        if (-not $shp.HasTextFrame) {
            continue
        }
        $tf = $shp.TextFrame
        if (-not $tf.HasText) {
            continue
        }
        $tr = $tf.TextRange
        $paraCount = $tr.Paragraphs().Count
        for ($i = 1; $i -le $paraCount; $i++) {
            $para = $tr.Paragraphs($i, 1)
            if ($para.IndentLevel -gt 1) {
                $para.IndentLevel = $para.IndentLevel - 1
            }
        }
    }
}
